$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new record as row 280 ("Región Metropolitana" / "Brócoli" /
# "Primera"), shifting all the existing rows 280-326 down to 281-327.
$ws.Rows.Item(280).Insert()

$ws.Range("A280").Value = 4
$ws.Range("B280").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C280").Value2 = "Los Lagos"
$ws.Range("D280").Value = 44694
$ws.Range("E280").Value = 10
$ws.Range("F280").Value = 100112023
$ws.Range("G280").Value2 = "Brócoli"
$ws.Range("H280").Value2 = "Sin especificar"
$ws.Range("I280").Value2 = "Primera"
$ws.Range("J280").Value = 1100
$ws.Range("K280").Value = 1500
$ws.Range("L280").Value = 1500
$ws.Range("M280").Value = 1500
$ws.Range("N280").Value2 = "$/unidad"
$ws.Range("O280").Value2 = "Región Metropolitana"
$ws.Range("P280").Value = 1500
$ws.Range("Q280").Value = 1
$ws.Range("R280").Value2 = "Hortaliza"
